$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain text (not parsed as a number by Excel),
# so a direct .Value assignment preserves the inline/shared string type.
$plainCells = @(
    @("B41", "Kaspa"),
    @("C41", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("E41", "  -0.10%  "),
    @("B42", "dogwifhat"),
    @("C42", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"),
    @("E42", "  +0.33%  "),
    @("D2", "63.408.76"),
    @("E2", "  +1.73%  "),
    @("D3", "3.172.77"),
    @("E3", "  -0.36%  "),
    @("E4", "  -0.12%  "),
    @("E5", "  +2.67%  "),
    @("E6", "  +0.71%  "),
    @("E7", "  -0.03%  "),
    @("D8", "3.170.69"),
    @("E8", "  -0.37%  "),
    @("E9", "  +2.16%  "),
    @("E10", "  +0.93%  "),
    @("E11", "  +2.46%  "),
    @("E12", "  +0.79%  "),
    @("E13", "  +2.15%  "),
    @("E14", "  +5.08%  "),
    @("D15", "3.695.37"),
    @("E15", "  -0.38%  "),
    @("E16", "  +1.37%  "),
    @("D17", "3.172.34"),
    @("E17", "  -0.44%  "),
    @("D18", "63.470.96"),
    @("E18", "  +1.73%  "),
    @("E19", "  +0.19%  "),
    @("E20", "  +1.04%  "),
    @("E21", "  +0.50%  "),
    @("E22", "  -0.79%  "),
    @("E23", "  +0.72%  "),
    @("E24", "  -0.24%  "),
    @("E25", "  +0.80%  "),
    @("E27", "  +0.84%  "),
    @("E28", "  -0.08%  "),
    @("E29", "  +3.83%  "),
    @("E30", "  -1.46%  "),
    @("E31", "  -1.38%  "),
    @("E32", "  -0.19%  "),
    @("E33", "  -2.19%  "),
    @("E34", "  +1.78%  "),
    @("E35", "  -1.49%  "),
    @("E36", "  +2.08%  "),
    @("D37", "0.0₃0736"),
    @("E37", "  +6.81%  "),
    @("E38", "  +0.30%  "),
    @("E39", "  +1.13%  "),
    @("E40", "  +1.59%  "),
    @("E43", "  -4.94%  "),
    @("D44", "2.800.85"),
    @("E44", "  -5.03%  "),
    @("E45", "  +0.47%  "),
    @("E46", "  +1.15%  "),
    @("E47", "  -1.20%  "),
    @("E49", "  +2.48%  "),
    @("E50", "  -1.46%  "),
    @("E51", "  +0.94%  ")
)

foreach ($pair in $plainCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Cells whose new value LOOKS like a number ("601.45", "0.999", ...).
# A direct .Value assignment would make Excel auto-convert the cell to a
# numeric type (matching the source data's original text representation,
# these must stay text). Route the literal through a scratch cell that is
# explicitly formatted as Text, then copy only the VALUE (paste-special)
# into the destination so the destination cell keeps its original (default)
# style while still being stored as a text string.
$textCells = @(
    @("D41", "0.113"),
    @("D42", "2.64"),
    @("D5", "601.45"),
    @("D6", "135.98"),
    @("D12", "0.454"),
    @("D14", "34.92"),
    @("D19", "6.59"),
    @("D20", "461.56"),
    @("D23", "7.66"),
    @("D24", "13.33"),
    @("D25", "83.24"),
    @("D26", "0.999"),
    @("D28", "1.00"),
    @("D30", "6.79"),
    @("D32", "27.18"),
    @("D33", "0.100"),
    @("D38", "51.30"),
    @("D39", "0.0391"),
    @("D40", "8.14"),
    @("D43", "391.48"),
    @("D46", "36.37"),
    @("D47", "2.12"),
    @("D49", "125.93"),
    @("D50", "25.15")
)

$scratch = $ws.Range("ZZ1")
foreach ($pair in $textCells) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $pair[1]
    $scratch.Copy()
    $ws.Range($pair[0]).PasteSpecial(-4163)
}
$scratch.Clear()

Write-Host "done"
